$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.190194845199585
$ws.Range("B1").Value = 2.088051080703735
$ws.Range("C1").Value = 4.258359909057617
$ws.Range("D1").Value = 3.001795530319214
$ws.Range("E1").Value = 1.213371396064758
